$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.10198523090934
$ws.Range("C2").Value = 7.488883920792741
$ws.Range("E2").Value = 11.55762992625142
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.654036658571503
$ws.Range("I2").Value = 23.97925797717541
$ws.Range("K2").Value = 9.420629954122841
$ws.Range("L2").Value = 9.781905511976325
$ws.Range("N2").Value = 19.23471626948237
$ws.Range("O2").Value = 23.98200986033059
$ws.Range("B3").Value = 11.86241284342279
$ws.Range("C3").Value = 7.47051833306339
$ws.Range("E3").Value = 11.56279376783696
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.655585473108621
$ws.Range("I3").Value = 24.07491580981009
$ws.Range("K3").Value = 9.260327322705969
$ws.Range("L3").Value = 9.768075586360188
$ws.Range("N3").Value = 19.29072635974914
$ws.Range("O3").Value = 24.07226700518478
$ws.Range("B4").Value = 11.71481957082065
$ws.Range("C4").Value = 7.459185677170469
$ws.Range("E4").Value = 11.56806108175095
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.656587227322757
$ws.Range("I4").Value = 24.1377168871862
$ws.Range("K4").Value = 9.161824376699574
$ws.Range("L4").Value = 9.7612382375489
$ws.Range("N4").Value = 19.32681670715059
$ws.Range("O4").Value = 24.13215381208198
$ws.Range("B5").Value = 11.65463240081571
$ws.Range("C5").Value = 7.45455445340172
$ws.Range("E5").Value = 11.57073553373033
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.657008256743487
$ws.Range("I5").Value = 24.16433160295987
$ws.Range("K5").Value = 9.121714442542686
$ws.Range("L5").Value = 9.758870459361017
$ws.Range("N5").Value = 19.34195244016689
$ws.Range("O5").Value = 24.15768094619072
$ws.Range("B6").Value = 11.64463830268926
$ws.Range("C6").Value = 7.453784682586273
$ws.Range("E6").Value = 11.57121153116125
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.657078942967969
$ws.Range("I6").Value = 24.16881273357428
$ws.Range("K6").Value = 9.115057555242359
$ws.Range("L6").Value = 9.758502635941126
$ws.Range("N6").Value = 19.34449163919788
$ws.Range("O6").Value = 24.16198749451809
$ws.Range("B7").Value = 11.71400792377589
$ws.Range("C7").Value = 7.459123270441986
$ws.Range("E7").Value = 11.56809501172683
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.656592853568193
$ws.Range("I7").Value = 24.13807168119877
$ws.Range("K7").Value = 9.161283248716526
$ws.Range("L7").Value = 9.761204607262354
$ws.Range("N7").Value = 19.32701909591334
$ws.Range("O7").Value = 24.13249353503419
$ws.Range("B8").Value = 12.01953236798755
$ws.Range("C8").Value = 7.482563733020182
$ws.Range("E8").Value = 11.55897575318466
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.654560174487234
$ws.Range("I8").Value = 24.01139682047815
$ws.Range("K8").Value = 9.365403274435755
$ws.Range("L8").Value = 9.776795194425421
$ws.Range("N8").Value = 19.25367643138215
$ws.Range("O8").Value = 24.0122025487876
$ws.Range("B9").Value = 12.61121946675051
$ws.Range("C9").Value = 7.528040611209659
$ws.Range("E9").Value = 11.55768949809942
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.650975223789751
$ws.Range("I9").Value = 23.79524937614853
$ws.Range("K9").Value = 9.762920780176978
$ws.Range("L9").Value = 9.820373031361376
$ws.Range("N9").Value = 19.12328613911651
$ws.Range("O9").Value = 23.81180737977305
$ws.Range("B10").Value = 13.03688270701272
$ws.Range("C10").Value = 7.561095390345357
$ws.Range("E10").Value = 11.56679859900699
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.64858343647336
$ws.Range("I10").Value = 23.65610421517572
$ws.Range("K10").Value = 10.05049482602392
$ws.Range("L10").Value = 9.860139613505236
$ws.Range("N10").Value = 19.03560236321828
$ws.Range("O10").Value = 23.6862689207723
$ws.Range("B11").Value = 13.22767662683619
$ws.Range("C11").Value = 7.57604274830233
$ws.Range("E11").Value = 11.57310799786879
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.647547391684781
$ws.Range("I11").Value = 23.5970708257089
$ws.Range("K11").Value = 10.179785293287
$ws.Range("L11").Value = 9.87986983984775
$ws.Range("N11").Value = 18.99745847358489
$ws.Range("O11").Value = 23.63387993383467
$ws.Range("B12").Value = 13.2994479826154
$ws.Range("C12").Value = 7.581688937451025
$ws.Range("E12").Value = 11.57580673900029
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.647162505272467
$ws.Range("I12").Value = 23.57532971474074
$ws.Range("K12").Value = 10.22848022981325
$ws.Range("L12").Value = 9.887572903645166
$ws.Range("N12").Value = 18.98326395923372
$ws.Range("O12").Value = 23.61472122150801
$ws.Range("B13").Value = 13.28401306378884
$ws.Range("C13").Value = 7.580473573578847
$ws.Range("E13").Value = 11.57521178396281
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.647245067050125
$ws.Range("I13").Value = 23.57998475847367
$ws.Range("K13").Value = 10.21800535980087
$ws.Range("L13").Value = 9.885903677750782
$ws.Range("N13").Value = 18.98630991035719
$ws.Range("O13").Value = 23.61881713601884
$ws.Range("B14").Value = 13.23359126517174
$ws.Range("C14").Value = 7.576507552138144
$ws.Range("E14").Value = 11.57332383358966
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.647515577913314
$ws.Range("I14").Value = 23.59526987291299
$ws.Range("K14").Value = 10.18379700581974
$ws.Range("L14").Value = 9.880498957812
$ws.Range("N14").Value = 18.99628568281429
$ws.Range("O14").Value = 23.63229010242851
$ws.Range("B15").Value = 13.20264219852848
$ws.Range("C15").Value = 7.574076384133996
$ws.Range("E15").Value = 11.57220765641564
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.647682241666311
$ws.Range("I15").Value = 23.60471235985143
$ws.Range("K15").Value = 10.16280765654301
$ws.Range("L15").Value = 9.877218450014189
$ws.Range("N15").Value = 19.00242862605508
$ws.Range("O15").Value = 23.64063125856359
$ws.Range("B16").Value = 13.02435064842523
$ws.Range("C16").Value = 7.560116644215166
$ws.Range("E16").Value = 11.56642968489458
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.648652187633942
$ws.Range("I16").Value = 23.66004802504771
$ws.Range("K16").Value = 10.04201075753208
$ws.Range("L16").Value = 9.858882848420398
$ws.Range("N16").Value = 19.03813016018283
$ws.Range("O16").Value = 23.68978773228058
$ws.Range("B17").Value = 12.91419582198806
$ws.Range("C17").Value = 7.551529046803204
$ws.Range("E17").Value = 11.56343846757429
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.649260509265596
$ws.Range("I17").Value = 23.69508711461312
$ws.Range("K17").Value = 9.967481770703424
$ws.Range("L17").Value = 9.848051737529783
$ws.Range("N17").Value = 19.06047781833442
$ws.Range("O17").Value = 23.72115320687019
$ws.Range("B18").Value = 12.85057402745505
$ws.Range("C18").Value = 7.546581342515525
$ws.Range("E18").Value = 11.56192196492459
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.649615295516424
$ws.Range("I18").Value = 23.71564203985253
$ws.Range("K18").Value = 9.924473352900183
$ws.Range("L18").Value = 9.841976621379766
$ws.Range("N18").Value = 19.073495804575
$ws.Range("O18").Value = 23.73963790593386
$ws.Range("B19").Value = 12.82898968253942
$ws.Range("C19").Value = 7.544904737674997
$ws.Range("E19").Value = 11.56144358912582
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.649736261964572
$ws.Range("I19").Value = 23.72267051008961
$ws.Range("K19").Value = 9.909888570142572
$ws.Range("L19").Value = 9.839946377037995
$ws.Range("N19").Value = 19.07793170464935
$ws.Range("O19").Value = 23.74597274486536
$ws.Range("B20").Value = 12.92594979451527
$ws.Range("C20").Value = 7.552444083529567
$ws.Range("E20").Value = 11.56373579112538
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.649195245953195
$ws.Range("I20").Value = 23.69131560136337
$ws.Range("K20").Value = 9.97543046546058
$ws.Range("L20").Value = 9.84918875148459
$ws.Range("N20").Value = 19.05808188487935
$ws.Range("O20").Value = 23.71776832592293
$ws.Range("B21").Value = 13.24841487915663
$ws.Range("C21").Value = 7.577672859650783
$ws.Range("E21").Value = 11.5738699866559
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.647435920639188
$ws.Range("I21").Value = 23.59076360985584
$ws.Range("K21").Value = 10.1938523545829
$ws.Range("L21").Value = 9.882080202781843
$ws.Range("N21").Value = 18.99334878659452
$ws.Range("O21").Value = 23.6283143043739
$ws.Range("B22").Value = 13.45634816273916
$ws.Range("C22").Value = 7.594078865026837
$ws.Range("E22").Value = 11.58229628696995
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.646329457191408
$ws.Range("I22").Value = 23.52862349965166
$ws.Range("K22").Value = 10.33504244343564
$ws.Range("L22").Value = 9.904924804261523
$ws.Range("N22").Value = 18.95249722719586
$ws.Range("O22").Value = 23.57381425586012
$ws.Range("B23").Value = 13.34564978669719
$ws.Range("C23").Value = 7.585330611996254
$ws.Range("E23").Value = 11.57763472080381
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.646916041847775
$ws.Range("I23").Value = 23.56146148497688
$ws.Range("K23").Value = 10.25984371012427
$ws.Range("L23").Value = 9.892610310161405
$ws.Range("N23").Value = 18.97416765087739
$ws.Range("O23").Value = 23.60253889764447
$ws.Range("B24").Value = 12.92063672870796
$ws.Range("C24").Value = 7.552030428299616
$ws.Range("E24").Value = 11.56360073810734
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.649224735751236
$ws.Range("I24").Value = 23.69301942388926
$ws.Range("K24").Value = 9.971837357701695
$ws.Range("L24").Value = 9.848674234218931
$ws.Range("N24").Value = 19.05916455671003
$ws.Range("O24").Value = 23.71929722189968
$ws.Range("B25").Value = 12.45242804478268
$ws.Range("C25").Value = 7.515796336961634
$ws.Range("E25").Value = 11.55626642703144
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.651902360141426
$ws.Range("I25").Value = 23.85027023445797
$ws.Range("K25").Value = 9.655972182684518
$ws.Range("L25").Value = 9.807209124916376
$ws.Range("N25").Value = 19.15712982529935
$ws.Range("O25").Value = 23.8622141122296
